# Apply the "created a nix file for running this analysis reproducably" edit.
#
# Summary of changes to data/output/tableau_upload.xlsx:
#  1. dc_data: append a 2024 row (District of Columbia, 2024, 1, 1, 1)
#  2. national_data: tweak a few aggregate counts (2021/2023/2024 rows)
#  3. national.monthly.strikes -> renamed to year.strikes.2024.monthly,
#     with several monthly values corrected and a new month-6 row added
#  4. a brand-new sheet year.strikes.2023.monthly is added with the full
#     12-month 2023 breakdown

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. dc_data — add the 2024 row
# ---------------------------------------------------------------------
$dc = $wb.Worksheets.Item("dc_data")
$dc.Cells.Item(4, 1).Value = "District of Columbia"
$dc.Cells.Item(4, 2).Value = 2024
$dc.Cells.Item(4, 3).Value = 1
$dc.Cells.Item(4, 4).Value = 1
$dc.Cells.Item(4, 5).Value = 1

# ---------------------------------------------------------------------
# 2. national_data — corrected totals
# ---------------------------------------------------------------------
$national = $wb.Worksheets.Item("national_data")
$national.Cells.Item(2, 2).Value = 90    # 2021 labor org count: 91 -> 90

$national.Cells.Item(4, 2).Value = 164   # 2023 labor org count: 166 -> 164
$national.Cells.Item(4, 4).Value = 473   # 2023 strikes: 475 -> 473

$national.Cells.Item(5, 2).Value = 86    # 2024 labor org count: 78 -> 86
$national.Cells.Item(5, 3).Value = 165   # 2024 employers: 144 -> 165
$national.Cells.Item(5, 4).Value = 230   # 2024 strikes: 193 -> 230

# ---------------------------------------------------------------------
# 3. national.monthly.strikes -> year.strikes.2024.monthly
# ---------------------------------------------------------------------
$monthly2024 = $wb.Worksheets.Item("national.monthly.strikes")
$monthly2024.Name = "year.strikes.2024.monthly"

# Month 2 (row 3): 14/16/17 -> 12/14/14
$monthly2024.Cells.Item(3, 2).Value = 12
$monthly2024.Cells.Item(3, 3).Value = 14
$monthly2024.Cells.Item(3, 4).Value = 14

# Month 3 (row 4): 22/47/48 -> 20/44/44
$monthly2024.Cells.Item(4, 2).Value = 20
$monthly2024.Cells.Item(4, 3).Value = 44
$monthly2024.Cells.Item(4, 4).Value = 44

# Month 4 (row 5): 20/28/29 -> 18/26/27
$monthly2024.Cells.Item(5, 2).Value = 18
$monthly2024.Cells.Item(5, 3).Value = 26
$monthly2024.Cells.Item(5, 4).Value = 27

# Month 5 (row 6): 12/14/14 -> 32/46/58
$monthly2024.Cells.Item(6, 2).Value = 32
$monthly2024.Cells.Item(6, 3).Value = 46
$monthly2024.Cells.Item(6, 4).Value = 58

# Month 6 (row 7): brand-new row
$monthly2024.Cells.Item(7, 1).Value = 6
$monthly2024.Cells.Item(7, 2).Value = 2
$monthly2024.Cells.Item(7, 3).Value = 2
$monthly2024.Cells.Item(7, 4).Value = 2

# ---------------------------------------------------------------------
# 4. new sheet: year.strikes.2023.monthly (full 12-month breakdown)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$monthly2023 = $wb.Worksheets.Add($null, $lastSheet)
$monthly2023.Name = "year.strikes.2023.monthly"

$monthly2023.Cells.Item(1, 1).Value = "Month"
$monthly2023.Cells.Item(1, 2).Value = "labor org count"
$monthly2023.Cells.Item(1, 3).Value = "employers"
$monthly2023.Cells.Item(1, 4).Value = "strikes"

$rows2023 = @(
    @(1, 40, 69, 95),
    @(2, 9, 10, 10),
    @(3, 19, 26, 28),
    @(4, 22, 31, 35),
    @(5, 16, 29, 40),
    @(6, 20, 27, 28),
    @(7, 15, 22, 31),
    @(8, 14, 25, 28),
    @(9, 23, 33, 48),
    @(10, 21, 32, 49),
    @(11, 24, 36, 39),
    @(12, 28, 39, 42)
)

$r = 2
foreach ($row in $rows2023) {
    $monthly2023.Cells.Item($r, 1).Value = $row[0]
    $monthly2023.Cells.Item($r, 2).Value = $row[1]
    $monthly2023.Cells.Item($r, 3).Value = $row[2]
    $monthly2023.Cells.Item($r, 4).Value = $row[3]
    $r++
}
